$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the shared "last_edited_time" timestamp text. Rows 4,5,6,7,8,12,13 in
# column D (last_edited_time) all held the same timestamp string; bump them all
# to the new value so they keep sharing a single string entry.
$oldTimestamp = "2024-08-03T21:27:00.000Z"
$newTimestamp = "2024-08-09T19:28:00.000Z"

foreach ($r in 4,5,6,7,8,12,13) {
    $cell = $ws.Cells.Item($r, 4)
    if ($cell.Value2 -eq $oldTimestamp) {
        $cell.Value2 = $newTimestamp
    }
}

# Update row 7 numeric metrics (properties.* .number / .formula.number columns).
$ws.Range("T7").Value2 = 10000000
$ws.Range("W7").Value2 = 42243000
$ws.Range("AA7").Value2 = 44507000
$ws.Range("AE7").Value2 = 86750000
$ws.Range("AH7").Value2 = 68250000
$ws.Range("AK7").Value2 = 11
$ws.Range("AN7").Value2 = 18500000
$ws.Range("AQ7").Value2 = 78250000
